$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.194.02"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.824.16"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.77"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5988"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06936"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2760"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.45"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07592"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "1.827.21"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.721"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6261"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009779"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.29"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "  -2.74%  "
$ws.Range("D17").Value = "29.007.11"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.521"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "  -7.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.70"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  -6.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.825"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  -3.23%  "
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "155.70"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.944"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1287"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.46"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06414"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "  -5.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.427"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.439"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.809"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.773"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").Value = "  -3.83%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6441"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.540"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.750"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01750"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.587"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").Value = "1.130.50"
$ws.Range("E40").Value = "  -8.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8882"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  -4.52%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "1.983.18"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.43"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.96"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("E46").Value = "  -3.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.610"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.391"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05493"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4531"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.336"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  -3.72%  "
